# "Big fix of ABS & RBA data sources"
#
# 1) watchlist: several C-column titles were stale/wrong. The ones that
#    previously showed a readable "Index Numbers ; ... ; Sydney/Melbourne/
#    Australia ;" label had actually failed to resolve and should show the
#    Error_<seriesId> placeholder instead; conversely the RBA/ABS non-market
#    sector rows (24, 34, 37) had been stuck showing the Error_<seriesId>
#    placeholder and should now show their resolved titles.
#
# 2) all_metadata: three new ABS series columns (A130272195C, A130272197J,
#    A130272199L) were appended after the existing last column (CG).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: watchlist sheet - fix up column C (title) for specific rows
# ---------------------------------------------------------------------
$wsWatch = $wb.Worksheets.Item("watchlist")

$watchFixes = @{
    2  = "Error_A2325941A"
    3  = "Error_A2326391L"
    4  = "Error_A2329946K"
    5  = "Error_A2326481T"
    6  = "Error_A2329856F"
    7  = "Error_A2328101R"
    8  = "Error_A2331566X"
    9  = "Error_A2331071W"
    10 = "Error_A2331386R"
    11 = "Error_A2329406W"
    12 = "Error_A3602793W"
    13 = "Error_A3602838R"
    14 = "Error_A2330531F"
    15 = "Error_A2330576K"
    16 = "Error_A2325896A"
    17 = "Error_A2325806K"
    21 = "Error_A2325811C"
    24 = "Gross value added non-market sector: Chain volume measures ;"
    34 = "Hours worked non-market sector: Index ;"
    37 = "Gross value added per hour worked non-market sector: Index ;"
    50 = "Error_A2325846C"
    51 = "Error_A2330616T"
    52 = "Error_A2330706W"
    53 = "Error_A2332236A"
}

foreach ($row in $watchFixes.Keys) {
    $wsWatch.Cells.Item($row, 3).Value = $watchFixes[$row]
}

# ---------------------------------------------------------------------
# Part 2: all_metadata sheet - append columns CH, CI, CJ (86, 87, 88)
#          for the three new series A130272195C / A130272197J / A130272199L
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("all_metadata")

$newIds = @("A130272195C", "A130272197J", "A130272199L")
$newCols = @(86, 87, 88)

for ($i = 0; $i -lt $newCols.Length; $i++) {
    $col = $newCols[$i]
    $id = $newIds[$i]

    $wsMeta.Cells.Item(1, $col).Value = $id              # row 1 - id header
    $wsMeta.Cells.Item(2, $col).Value = "Unknown"         # row 2 - units
    $wsMeta.Cells.Item(5, $col).Value = "Unknown"         # row 5 - frequency
    $wsMeta.Cells.Item(8, $col).Value = "abs_series"      # row 8 - source
    $wsMeta.Cells.Item(17, $col).Value = $id              # row 17 - title
    $wsMeta.Cells.Item(18, $col).Value = $id              # row 18 - id
}
